$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking
# strings (e.g. "216.07") are not auto-converted to numbers, matching
# the inline-string cell type used in the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.071.05"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "1.668.61"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "216.07"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D8").Value = "0.2687"
$ws.Range("D9").Value = "0.06372"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").Value = "21.84"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").Value = "1.674.90"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "4.512"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "0.5810"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "0.000008492"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").Value = "64.13"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "25.899.86"
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").Value = "4.927"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "10.81"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "189.74"
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("D22").Value = "6.190"
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "144.61"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").Value = "7.602"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").Value = "0.1221"
$ws.Range("E26").Value = "  +4.25%  "
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").Value = "0.06647"
$ws.Range("E28").Value = "  +13.53%  "
$ws.Range("D29").Value = "1.340"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("D32").Value = "3.524"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("D33").Value = "1.664"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").Value = "0.6168"
$ws.Range("E35").Value = "  +3.65%  "
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").Value = "2.685"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").Value = "6.247"
$ws.Range("E38").Value = "  +5.97%  "
$ws.Range("D39").Value = "1.095.34"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").Value = "0.01596"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").Value = "0.8678"
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("D44").Value = "1.815.39"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("D45").Value = "0.00000000115"
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("D46").Value = "56.38"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "1.004"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "8.124"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").Value = "0.05235"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("D50").Value = "0.4281"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "5.992"
$ws.Range("E51").Value = "  +2.58%  "
